# EPBDS-6561 Add cycled dependency on beans
#
# Adds a new field (Address adr) to the existing "Datatype Person" table
# (row 7), and a new "Environment" / "import com.example.beans" table
# (rows 10-11) to the datatypes sheet, matching the diff against
# xl/worksheets/sheet1.xml and xl/sharedStrings.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: new "Address adr" field appended to the Person datatype ---
$ws.Range("B7").Value = "Address"
$ws.Range("C7").Value = "adr"

# --- Rows 10-11: new Environment / import table -----------------------
$ws.Range("B10").Value = "Environment"
# C10 stays empty but must still exist as a (styled) cell in the sheet.
$ws.Range("C10").Style = "Normal"

$ws.Range("B11").Value = "import"
$ws.Range("C11").Value = "com.example.beans"

# --- Column widths (bestFit in the original) ---------------------------
$ws.Columns("B").ColumnWidth = 14.833333333333334
$ws.Columns("C").ColumnWidth = 18.333333333333332
$ws.Columns("D").ColumnWidth = 24.166666666666668

# --- Selection ends up on C10, matching the recorded sheetView ---------
$ws.Range("C10").Select()
